# Add results from Kansas City 2017.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Run time, ms" results for the Kansas City 2017 run, entered into column D
# alongside the existing "Sample run time, ms" values in column E.
$ws.Range("D3").Value = 15306
$ws.Range("D4").Value = 10421
$ws.Range("D5").Value = 9200
$ws.Range("D6").Value = 131

$ws.Range("D8").Value = 3479
$ws.Range("D9").Value = 3113
$ws.Range("D10").Value = 926
$ws.Range("D11").Value = 467

$ws.Range("D13").Value = 3151
$ws.Range("D14").Value = 3575
$ws.Range("D15").Value = 1416
$ws.Range("D16").Value = 230

$ws.Range("D18").Value = 183561
$ws.Range("D19").Value = 116
$ws.Range("D20").Value = 194591
$ws.Range("D21").Value = 184538
$ws.Range("D22").Value = 1901
$ws.Range("D23").Value = 1779
$ws.Range("D24").Value = 249

$ws.Range("D26").Value = 5158
$ws.Range("D27").Value = 100000000
$ws.Range("D28").Value = 555

$ws.Range("D30").Value = 2279
$ws.Range("D31").Value = 565

$ws.Range("D33").Value = 2739
$ws.Range("D34").Value = 85

# Move the on-screen selection/scroll position to reflect where the user
# finished entering data (just past the last data row).
$ws.Range("D35").Select()
$excel.ActiveWindow.ScrollRow = 23
